# Calibration Legs and Others Update!
# Apply the updated calibration readings for rows 5, 6, 8 and 9. All the
# formula cells (H, I, J, K, L, M) are derived from these inputs and will
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (LEG L2)
$ws.Range("D5").Value = 1600
$ws.Range("E5").Value = 1000

# Row 6 (LEG L3)
$ws.Range("B6").Value = 1400

# Row 8 (LEG R2)
$ws.Range("B8").Value = 1400
$ws.Range("C8").Value = 1500
$ws.Range("D8").Value = 1200
$ws.Range("E8").Value = 1850

# Row 9 (LEG R3)
$ws.Range("B9").Value = 1300
$ws.Range("E9").Value = 1800

# Scroll the sheet view so column E is at the left edge (best effort — mirrors
# the authored view state `topLeftCell="E1"`).
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
